$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "向 ROW 客戶傳送提醒電子郵件" "向 ROW 客户发送提醒电子邮件"
Replace-Text "將於 9 月 29 日移除 Tether Omni (USDT)" "将于 9 月 29 日移除 Tether Omni (USDT)"
Replace-Text "向 Tether Omni 道別" "向 Tether Omni 说再见"
Replace-Text "自 2023 年 9 月 29 日格林威治標準時間 00:00 起，Deriv 將停止提供Tether Omni (USDT) 作為帳戶貨幣。 這是因為 Tether 已停止支援 USDT 的 Omni 轉帳。" "自 2023 年 9 月 29 日格林威治标准时间 00:00 起，Deriv 将停止提供Tether Omni (USDT) 作为账户货币。 这是因为 Tether 已停止支持 USDT 的 Omni 转账。"
Replace-Text "需要做什麼？" "需要做什么？"
Replace-Text "如果 USDT 帳戶 " "如果 USDT 账户 "
Replace-Text " 中有餘額，請在上述日期之前提取餘額。 若有持倉頭寸，提取餘額前請先平倉。" " 中有余额，请在上述日期之前提取余额。 如果有持仓头寸，提取余额之前请先平仓。"
Replace-Text "檢查帳戶" "查看账户"
Replace-Text "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何持倉頭寸將在上述日期後自動平倉，帳戶餘額將轉移到最後活躍的帳戶" "USDT 账户将于 2023 年 9 月 29 日格林尼治标准时间 00:00 关闭。 任何持仓头寸将在上述日期后自动平仓，账户余额将转移到最后活跃的账户"
Replace-Text "在此過程中將採用標準匯率和費用。" "在此过程中将采用标准汇率和费用。"
Replace-Text "如有任何疑問，請透過以下方式聯繫我們：" "如有任何疑问，请通过以下方式联系我们："
Replace-Text "即時聊天" "实时聊天"

foreach ($c in $d.Comments) {
    $c.Range.Find.Execute("@azita@regentmarkets.com，BE 無法保證屆時可以完成指令，", $true, $false, $false, $false, $false, $true, 1, $false, "@azita@regentmarkets.com，BE 不能保证届时可以完成脚本，", 2) | Out-Null
    $c.Range.Find.Execute("可以將其稱為“所述日期之後”嗎？", $true, $false, $false, $false, $false, $true, 1, $false, "可以将其称为“所述日期之后”吗？", 2) | Out-Null
    $c.Range.Find.Execute("意思是在所述日期之後完成轉帳嗎？", $true, $false, $false, $false, $false, $true, 1, $false, "意思是在所述日期之后完成转账吗？", 2) | Out-Null
    $c.Range.Find.Execute("是的... 看來我們無法確認日期", $true, $false, $false, $false, $false, $true, 1, $false, "是的... 看来我们无法确认日期", 2) | Out-Null
}
